# realization: WR &beam dynamic rework done
#
# The "Flat-top >= length of batch" caption (present once on the "rise /
# flat-top / fall" slide and once again on the injection-path slide that
# reuses the same artwork) is reworded to call out the *kicker* flat-top
# explicitly, and is wrapped onto two lines. Because the caption now spans
# two lines, its text box is grown taller and nudged so it stays visually
# anchored to the same arrow.
#
# Note: this COM host only reliably binds *positional* parameters on
# user-defined functions, so the helpers below avoid named arguments.

$p = $ppt.ActivePresentation

function Get-ShapeById {
    param($Slide, $Id)
    for ($i = 1; $i -le $Slide.Shapes.Count; $i++) {
        $candidate = $Slide.Shapes.Item($i)
        if ($candidate.Id -eq $Id) {
            return $candidate
        }
    }
    return $null
}

function Update-FlatTopCaption {
    param($Slide, $Left, $Top, $Width, $Height)

    $shape = Get-ShapeById $Slide 94

    # Split the single line into two paragraphs. Re-assigning .Text like
    # this preserves the existing run formatting (font/size/panose/etc.)
    # on each resulting paragraph, since neither run's properties are
    # touched directly.
    $shape.TextFrame.TextRange.Text = "kicker flat-top ≥ " + [char]13 + "length of batch"

    # Re-anchor/resize the box now that it holds two lines of text. Values
    # are expressed in points (as the Shape.Left/Top/Width/Height COM
    # properties require) but carry extra precision so that, once rounded
    # to the single-precision float the property actually stores, they
    # still convert back to the exact target EMU coordinates.
    $shape.Left = $Left
    $shape.Top = $Top
    $shape.Width = $Width
    $shape.Height = $Height
}

# --- Slide 1 ("rise / flat-top / fall") ------------------------------------
$slide1 = $p.Slides.Item(1)
Update-FlatTopCaption $slide1 390.3285217285 596.2550048828 290.4769592285 64.8000030518

# --- Slide 3 (beam-injection path reuses the same caption) -----------------
$slide3 = $p.Slides.Item(3)
Update-FlatTopCaption $slide3 442.1785278320 640.5050048828 290.4769592285 64.8000030518
